$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (Employee_Name, Employee_No, Dept_Name, Dept_Number, Total_Compensation, Months_Spent)
$data = @(
    @("WARD",   "7521", "SALES",      "30", "623750",  "499"),
    @("TURNER", "7844", "SALES",      "30", "739500",  "493"),
    @("SCOTT",  "7788", "RESEARCH",   "20", "1184480", "410"),
    @("JAMES",  "7900", "SALES",      "30", "452200",  "476"),
    @("MILLER", "7934", "ACCOUNTING", "10", "634400",  "488"),
    @("SMITH",  "7369", "RESEARCH",   "20", "400800",  "501"),
    @("MARTIN", "7654", "SALES",      "30", "615000",  "492"),
    @("CLARK",  "7782", "ACCOUNTING", "10", "1215200", "496"),
    @("JAMES",  "7900", "ACCOUNTING", "30", "12350",   "13"),
    @("JONES",  "7566", "RESEARCH",   "20", "1481550", "498"),
    @("FORD",   "7902", "RESEARCH",   "20", "1470000", "490"),
    @("ALLEN",  "7499", "SALES",      "30", "798400",  "499"),
    @("BLAKE",  "7698", "SALES",      "30", "1416450", "497"),
    @("KING",   "7839", "ACCOUNTING", "10", "2450000", "490"),
    @("ADAMS",  "7876", "RESEARCH",   "20", "465300",  "423")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}

# Remove the now-unused rows 17 and 18 (shrinks the range from A1:F18 to A1:F16)
$ws.Range("A17:F18").Delete()
